$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New feedback entry (row 55) mirrors the layout/formatting of the row above it.
$ws.Range("A54:F54").Copy($ws.Range("A55:F55"))

$ws.Range("A55").Value = 41929
$ws.Range("B55").Value = "Email from SAM"
$ws.Range("C55").Value = "Cook, Jeff <Jeff.Cook@nrel.gov>"
$ws.Range("D55").Value = "I am a little curious as to why you can only adjust the rows and number of turbines per row, to get the total number of turbines to change. To me it would be more intuitive if you could change the number of turbines in the first cell, and have that adjust your rows and number of turbines per row by default. "
$ws.Range("E55").Value = "Replied with cc to Janine. Should be addressed by new Wind Farm layout option"
$ws.Range("F55").Value = 41929

$ws.Rows.Item(55).RowHeight = 75

# Move the selection to reflect where editing left off.
$ws.Range("E56").Select()
